$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 154
$ws.Range("I4").Value = 138.36363
$ws.Range("J4").Value = 240
$ws.Range("K4").Value = 138.36363
$ws.Range("L4").Value = 240
$ws.Range("M4").Value = -24.36363
$ws.Range("N4").Value = -468

$ws.Range("H17").Value = 2386763.5
$ws.Range("J17").Value = 2386763.5
$ws.Range("L17").Value = 7160290.5
$ws.Range("N17").Value = -7160626.5

$ws.Range("H40").Value = 2171.4827
$ws.Range("I40").Value = 1988.421
$ws.Range("J40").Value = 2519.3
$ws.Range("K40").Value = 1988.421
$ws.Range("L40").Value = 2519.3
$ws.Range("M40").Value = -1813.421
$ws.Range("N40").Value = -2869.3

$ws.Range("H112").Value = 1115.2449
$ws.Range("I112").Value = 1163.3334
$ws.Range("J112").Value = 1112.1086
$ws.Range("K112").Value = 3490.0002
$ws.Range("L112").Value = 3336.3258
$ws.Range("M112").Value = -2382.0002
$ws.Range("N112").Value = -5552.325800000001

$ws.Range("H132").Value = 3239.0688
$ws.Range("I132").Value = 2831.913
$ws.Range("J132").Value = 4799.8335
$ws.Range("K132").Value = 8495.739
$ws.Range("L132").Value = 14399.5005
$ws.Range("M132").Value = -5965.739
$ws.Range("N132").Value = -19459.5005

$ws.Range("H138").Value = 2467.4736
$ws.Range("I138").Value = 1529.7106
$ws.Range("J138").Value = 3092.6492
$ws.Range("K138").Value = 4589.1318
$ws.Range("L138").Value = 9277.9476
$ws.Range("M138").Value = 550.8681999999999
$ws.Range("N138").Value = -19557.9476

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6687.2656
$ws.Range("I32").Value = 5656.1577
$ws.Range("J32").Value = 15083.429
$ws.Range("K32").Value = 5656.1577
$ws.Range("L32").Value = 15083.429
$ws.Range("M32").Value = -5369.1577
$ws.Range("N32").Value = -15657.429

$ws.Range("H45").Value = 1954.1538
$ws.Range("I45").Value = 1967.1111
$ws.Range("K45").Value = 1967.1111
$ws.Range("M45").Value = -1590.1111

$ws.Range("H74").Value = 14402.6
$ws.Range("I74").Value = 110012
$ws.Range("J74").Value = 3779.3333
$ws.Range("K74").Value = 110012
$ws.Range("L74").Value = 3779.3333
$ws.Range("M74").Value = -109138
$ws.Range("N74").Value = -5527.3333

$ws.Range("H77").Value = 14402.6
$ws.Range("I77").Value = 110012
$ws.Range("J77").Value = 3779.3333
$ws.Range("K77").Value = 550060
$ws.Range("L77").Value = 18896.6665
$ws.Range("M77").Value = -545692
$ws.Range("N77").Value = -27632.6665

$ws.Range("H122").Value = 4338
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4338
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 13014
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -17914

$ws.Range("H132").Value = 2516.1538
$ws.Range("I132").Value = 1809.8
$ws.Range("J132").Value = 4870.6665
$ws.Range("K132").Value = 5429.4
$ws.Range("L132").Value = 14611.9995
$ws.Range("M132").Value = -2899.4
$ws.Range("N132").Value = -19671.9995

$ws.Range("H139").Value = 38942.273
$ws.Range("J139").Value = 38942.273
$ws.Range("L139").Value = 38942.273
$ws.Range("N139").Value = -49222.273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 483
$ws.Range("I22").Value = 474.5
$ws.Range("K22").Value = 474.5
$ws.Range("M22").Value = -301.5

$ws.Range("H86").Value = 4305.2983
$ws.Range("I86").Value = 4397.244
$ws.Range("J86").Value = 4069.6875
$ws.Range("K86").Value = 4397.244
$ws.Range("L86").Value = 4069.6875
$ws.Range("M86").Value = -3274.244
$ws.Range("N86").Value = -6315.6875

$ws.Range("H89").Value = 4305.2983
$ws.Range("I89").Value = 4397.244
$ws.Range("J89").Value = 4069.6875
$ws.Range("K89").Value = 21986.22
$ws.Range("L89").Value = 20348.4375
$ws.Range("M89").Value = -16370.22
$ws.Range("N89").Value = -31580.4375

$ws.Range("H105").Value = 2628.3845
$ws.Range("I105").Value = 2558.25
$ws.Range("J105").Value = 3470
$ws.Range("K105").Value = 2558.25
$ws.Range("L105").Value = 3470
$ws.Range("M105").Value = -811.25
$ws.Range("N105").Value = -6964

$ws.Range("H107").Value = 876.8461
$ws.Range("I107").Value = 839.8889
$ws.Range("J107").Value = 960
$ws.Range("K107").Value = 839.8889
$ws.Range("L107").Value = 960
$ws.Range("M107").Value = 1080.1111
$ws.Range("N107").Value = -4800

$ws.Range("H134").Value = 1563.4878
$ws.Range("I134").Value = 1283.303
$ws.Range("J134").Value = 2719.25
$ws.Range("K134").Value = 3849.909000000001
$ws.Range("L134").Value = 8157.75
$ws.Range("M134").Value = -1314.909000000001
$ws.Range("N134").Value = -13227.75

$ws.Range("H140").Value = 44610
$ws.Range("J140").Value = 44610
$ws.Range("L140").Value = 44610
$ws.Range("N140").Value = -54970

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2207.5
$ws.Range("I58").Value = 1651.3
$ws.Range("J58").Value = 3002.0715
$ws.Range("K58").Value = 1651.3
$ws.Range("L58").Value = 3002.0715
$ws.Range("M58").Value = -1448.3
$ws.Range("N58").Value = -3408.0715

$ws.Range("H122").Value = 455958.9
$ws.Range("I122").Value = 998.2353000000001
$ws.Range("J122").Value = 2002825.2
$ws.Range("K122").Value = 2994.7059
$ws.Range("L122").Value = 6008475.6
$ws.Range("M122").Value = -544.7058999999999
$ws.Range("N122").Value = -6013375.6

$ws.Range("H136").Value = 2207.5
$ws.Range("I136").Value = 1651.3
$ws.Range("J136").Value = 3002.0715
$ws.Range("K136").Value = 4953.9
$ws.Range("L136").Value = 9006.2145
$ws.Range("M136").Value = -2403.9
$ws.Range("N136").Value = -14106.2145

$ws.Range("H138").Value = 34836.11
$ws.Range("J138").Value = 34836.11
$ws.Range("L138").Value = 34836.11
$ws.Range("N138").Value = -45116.11

$ws.Range("H140").Value = 56595.2
$ws.Range("J140").Value = 56595.2
$ws.Range("L140").Value = 56595.2
$ws.Range("N140").Value = -66955.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 892.29
$ws.Range("J131").Value = 897.4388
$ws.Range("L131").Value = 2692.3164
$ws.Range("N131").Value = -12772.3164

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4060.2856
$ws.Range("I132").Value = 3684.8
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 11054.4
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -8524.400000000001
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 19201
$ws.Range("I40").Value = 26371.7
$ws.Range("J40").Value = 3266.111
$ws.Range("K40").Value = 26371.7
$ws.Range("L40").Value = 3266.111
$ws.Range("M40").Value = -26235.7
$ws.Range("N40").Value = -3538.111

$ws.Range("H61").Value = 548.63635
$ws.Range("I61").Value = 588.6111
$ws.Range("J61").Value = 368.75
$ws.Range("K61").Value = 588.6111
$ws.Range("L61").Value = 368.75
$ws.Range("M61").Value = -386.6111
$ws.Range("N61").Value = -772.75

$ws.Range("H113").Value = 548.63635
$ws.Range("I113").Value = 588.6111
$ws.Range("J113").Value = 368.75
$ws.Range("K113").Value = 588.6111
$ws.Range("L113").Value = 368.75
$ws.Range("M113").Value = 1581.3889
$ws.Range("N113").Value = -4708.75

$ws.Range("H127").Value = 52340.668
$ws.Range("J127").Value = 52340.668
$ws.Range("L127").Value = 52340.668
$ws.Range("N127").Value = -62260.668

$ws.Range("H132").Value = 14133.936
$ws.Range("I132").Value = 12140.477
$ws.Range("J132").Value = 18320.2
$ws.Range("K132").Value = 36421.431
$ws.Range("L132").Value = 54960.60000000001
$ws.Range("M132").Value = -33891.431
$ws.Range("N132").Value = -60020.60000000001

$ws.Range("H136").Value = 13376109
$ws.Range("I136").Value = 58104.223
$ws.Range("K136").Value = 174312.669
$ws.Range("M136").Value = -171762.669

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 279290.47
$ws.Range("I132").Value = 385709.56
$ws.Range("J132").Value = 2600.8
$ws.Range("K132").Value = 1157128.68
$ws.Range("L132").Value = 7802.400000000001
$ws.Range("M132").Value = -1154598.68
$ws.Range("N132").Value = -12862.4

$ws.Range("H136").Value = 306643.06
$ws.Range("I136").Value = 455397.3
$ws.Range("K136").Value = 1366191.9
$ws.Range("M136").Value = -1363641.9
